$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.795.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.290.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -1.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.633.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.852'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.282.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.789.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000111'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.00%  '
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +11.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.36%  '
$ws.Range("E25").Value = '  +7.21%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("E35").Value = '  +0.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.68'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0351'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.237'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.93%  '
$ws.Range("E41").Value = '  -1.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.39'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.80%  '
$ws.Range("E45").Value = '  -2.74%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.103'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.445'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.56%  '
